# Apply LDLC price-history update:
# Insert a new column at GX (shifting old GX->GY, old GY->GZ).
# New GX column gets a fresh timestamp header and, for rows that
# currently carry a numeric price in GW (rows 2-80), the same price
# value (mirroring the latest price check). Rows without a price
# (81-210) are left blank, matching the surrounding empty cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before GX; this shifts GX -> GY and GY -> GZ,
# leaving the new GX column empty.
$ws.Range("GX:GX").EntireColumn.Insert()

# New header cell for the freshly inserted column, matching the bold /
# bordered / centered style used by the other header cells.
$ws.Range("GX1").Value = "2026-02-06 13:53:41"
$ws.Range("GX1").Font.Bold = $true
$ws.Range("GX1").HorizontalAlignment = -4108
$ws.Range("GX1").VerticalAlignment = -4160
$ws.Range("GX1").Borders.LineStyle = 1

# Find the last used row.
$lastRow = $ws.UsedRange.Rows.Count()

# Mirror the price from GW into the new GX column wherever GW holds
# a numeric price (rows 2-80 in the original data).
for ($r = 2; $r -le $lastRow; $r++) {
    $gw = $ws.Cells.Item($r, 205)
    $val = $gw.Value()
    if ($val -ne $null -and $val -ne "") {
        $gx = $ws.Cells.Item($r, 206)
        $gx.Value = $val
    }
}
